$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.837392449378967
$ws.Range("B1").Value = 2.079525470733643
$ws.Range("C1").Value = 2.464934110641479
$ws.Range("D1").Value = 2.417133808135986
$ws.Range("E1").Value = 2.570950031280518
